$d = $word.ActiveDocument

# --- Edit 1: first paragraph, append "  (This is a change – Version for main branch)" ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$textLen = $r.Text.Length
$end = $r.Start + $textLen - 1   # position right before the paragraph mark

# two trailing spaces on the plain (black) run
$ip = $d.Range($end, $end)
$ip.InsertAfter("  ")
$end = $end + 2

$seg1 = "(This is a change " + [char]0x2013 + " Ve"
$ip = $d.Range($end, $end)
$ip.InsertAfter($seg1)
$seg1Range = $d.Range($end, $end + $seg1.Length)
$seg1Range.Font.Color = 255
$end = $end + $seg1.Length

$seg2 = "rsion for main branch"
$ip = $d.Range($end, $end)
$ip.InsertAfter($seg2)
$seg2Range = $d.Range($end, $end + $seg2.Length)
$seg2Range.Font.Color = 255
$end = $end + $seg2.Length

$seg3 = ")"
$ip = $d.Range($end, $end)
$ip.InsertAfter($seg3)
$seg3Range = $d.Range($end, $end + $seg3.Length)
$seg3Range.Font.Color = 255
$end = $end + $seg3.Length

Write-Output ("Paragraph1: [" + $p1.Range.Text + "]")

# --- Edit 2: add a new empty paragraph at the end of the document with shading F9F9F9 ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endOfDoc = $lastPara.Range.End
$tail = $d.Range($endOfDoc, $endOfDoc)
$tail.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Shading.BackgroundPatternColor = 16382457
Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
